$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing of existing answers
$ws.Range("B2").Value = "nick"
$ws.Range("B3").Value = "blue"

# Add new question/answer rows
$ws.Range("A5").Value = "sadf?"
$ws.Range("B5").Value = "bob"
$ws.Range("A6").Value = "who eat cats?"
$ws.Range("B6").Value = "me"
$ws.Range("A7").Value = "who farted?"
$ws.Range("B7").Value = "me"

# Update the selection to match the saved state
$ws.Range("C3").Select()
